$d = $word.ActiveDocument
$p51 = $d.Paragraphs.Item(51)
$p53 = $d.Paragraphs.Item(53)
$startPos = $p51.Range.Start
$endPos = $p53.Range.End
$rngDelete = $d.Range($startPos, $endPos)
$rngDelete.Delete()

$p50 = $d.Paragraphs.Item(50)
$full = $p50.Range
$xml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p w14:paraId="140E3111" w14:textId="0CEB3880" w:rsidR="00966A4F" w:rsidRPr="003A07A8" w:rsidRDefault="00966A4F" w:rsidP="00966A4F"><w:pPr><w:pStyle w:val="NormalWeb"/><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r w:rsidRPr="00966A4F"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr><w:t xml:space="preserve">In this case, it is greater, so </w:t></w:r><w:r w:rsidR="003A07A8" w:rsidRPr="003A07A8"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/></w:rPr><w:t>R</w:t></w:r><w:r w:rsidRPr="003A07A8"><w:rPr><w:rFonts w:ascii="Calibri" w:eastAsiaTheme="minorHAnsi" w:hAnsi="Calibri" w:cs="Calibri"/><w:b/></w:rPr><w:t>eject the null.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$full.InsertXML($xml)
Write-Host "done"
